# "aggiornamento fino a 28/06 incluso"
# Append daily COVID-style rows 270-301 (dates 2021-05-28 .. 2021-06-28,
# serials 44344..44375) to the existing data table that currently ends at
# row 269 (serial 44343 / 2021-05-27).
#
# Columns: A = date serial, B = nuovi pos., C = somma mobile 7gg.,
#          D = somma mobile 7gg. per 100mila abitanti.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(44344, 0, 1, 31.25976867771178),
    @(44345, 1, 2, 62.51953735542357),
    @(44346, 0, 2, 62.51953735542357),
    @(44347, 0, 2, 62.51953735542357),
    @(44348, 0, 2, 62.51953735542357),
    @(44349, 0, 2, 62.51953735542357),
    @(44350, 0, 1, 31.25976867771178),
    @(44351, 1, 2, 62.51953735542357),
    @(44352, 0, 1, 31.25976867771178),
    @(44353, 0, 1, 31.25976867771178),
    @(44354, 0, 1, 31.25976867771178),
    @(44355, 0, 1, 31.25976867771178),
    @(44356, 0, 1, 31.25976867771178),
    @(44357, 0, 1, 31.25976867771178),
    @(44358, 0, 0, 0),
    @(44359, 0, 0, 0),
    @(44360, 0, 0, 0),
    @(44361, 0, 0, 0),
    @(44362, 0, 0, 0),
    @(44363, 0, 0, 0),
    @(44364, 0, 0, 0),
    @(44365, 0, 0, 0),
    @(44366, 0, 0, 0),
    @(44367, 0, 0, 0),
    @(44368, 0, 0, 0),
    @(44369, 0, 0, 0),
    @(44370, 0, 0, 0),
    @(44371, 0, 0, 0),
    @(44372, 0, 0, 0),
    @(44373, 0, 0, 0),
    @(44374, 0, 0, 0),
    @(44375, 0, 0, 0)
)

$row = 270
foreach ($entry in $newRows) {
    $prevRow = $row - 1

    # Copy the formatting of the row directly above (style "s=2" on the date
    # cell in column A, etc.) down into the new row before writing values,
    # so the appended rows keep the same look as the rest of the table.
    $ws.Range("A" + $prevRow + ":D" + $prevRow).Copy()
    $ws.Range("A" + $row + ":D" + $row).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]

    $row = $row + 1
}

$excel.CutCopyMode = 0
